$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert textual TRUE/FALSE answers and numeric short-answer values
# to their proper native Excel types.
$ws.Range("G3").Value = $true
$ws.Range("G4").Value = 4
$ws.Range("G6").Value = $true

# Update the last active selection/cursor position to match the authored file.
$ws.Range("D18").Select()
